# Update the two-digit division answers table.
# The table has 20 rows x 5 columns, but only every 4th row (1,5,9,13,17)
# contains the answer text; the remaining rows are blank spacer rows.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "45÷4=11, 1" },
    @{ Row = 1;  Col = 2; New = "81÷8=10, 1" },
    @{ Row = 1;  Col = 3; New = "42÷3=14, 0" },
    @{ Row = 1;  Col = 4; New = "47÷4=11, 3" },
    @{ Row = 1;  Col = 5; New = "37÷8=4, 5" },

    @{ Row = 5;  Col = 1; New = "85÷5=17, 0" },
    @{ Row = 5;  Col = 2; New = "72÷2=36, 0" },
    @{ Row = 5;  Col = 3; New = "44÷3=14, 2" },
    @{ Row = 5;  Col = 4; New = "29÷7=4, 1" },
    @{ Row = 5;  Col = 5; New = "80÷7=11, 3" },

    @{ Row = 9;  Col = 1; New = "76÷8=9, 4" },
    @{ Row = 9;  Col = 2; New = "17÷4=4, 1" },
    @{ Row = 9;  Col = 3; New = "44÷9=4, 8" },
    @{ Row = 9;  Col = 4; New = "28÷9=3, 1" },
    @{ Row = 9;  Col = 5; New = "82÷3=27, 1" },

    @{ Row = 13; Col = 1; New = "61÷7=8, 5" },
    @{ Row = 13; Col = 2; New = "85÷5=17, 0" },
    @{ Row = 13; Col = 3; New = "21÷2=10, 1" },
    @{ Row = 13; Col = 4; New = "70÷4=17, 2" },
    @{ Row = 13; Col = 5; New = "87÷3=29, 0" },

    @{ Row = 17; Col = 1; New = "94÷3=31, 1" },
    @{ Row = 17; Col = 2; New = "18÷7=2, 4" },
    @{ Row = 17; Col = 3; New = "13÷7=1, 6" },
    @{ Row = 17; Col = 4; New = "65÷8=8, 1" },
    @{ Row = 17; Col = 5; New = "69÷5=13, 4" }
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.New
}
